$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.94
$ws.Range("C4").Value = -11.513
$ws.Range("B7").Value = 5.553000000000001
$ws.Range("A8").Value = -22.191
$ws.Range("A10").Value = -21.317
$ws.Range("D10").Value = -8.164000000000001
$ws.Range("C11").Value = -12.126
$ws.Range("A12").Value = -21.618
$ws.Range("D12").Value = -7.152000000000001
$ws.Range("D13").Value = -8.054
$ws.Range("B14").Value = 5.614000000000001
$ws.Range("C14").Value = -12.108
$ws.Range("D14").Value = -7.481999999999999
$ws.Range("B15").Value = 5.242
$ws.Range("A18").Value = -21.529
$ws.Range("B18").Value = 6.498
$ws.Range("C18").Value = -11.806
$ws.Range("C19").Value = -11.752
$ws.Range("B20").Value = 6.489999999999999
$ws.Range("C21").Value = -12.083
$ws.Range("A25").Value = -21.789
$ws.Range("C27").Value = -13.469
$ws.Range("B29").Value = 4.892999999999999
$ws.Range("D29").Value = -7.170999999999999
$ws.Range("B30").Value = 5.205
$ws.Range("B31").Value = 5.049
$ws.Range("C31").Value = -13.286
$ws.Range("D32").Value = -8.296000000000001
$ws.Range("B35").Value = 8.379
$ws.Range("D35").Value = -7.781000000000001
$ws.Range("A37").Value = -20.413
$ws.Range("C38").Value = -13.052
$ws.Range("B40").Value = 8.382000000000001
$ws.Range("C42").Value = -12.533
$ws.Range("D43").Value = -8.657999999999998
$ws.Range("B44").Value = 5.577
$ws.Range("C44").Value = -13.28
$ws.Range("C47").Value = -12.241
$ws.Range("D48").Value = -7.461
$ws.Range("D49").Value = -8.144
$ws.Range("B50").Value = 5.210999999999999
$ws.Range("D50").Value = -8.006
$ws.Range("D51").Value = -8.041
$ws.Range("B54").Value = 5.066999999999999
$ws.Range("A55").Value = -21.837
$ws.Range("C56").Value = -12.771
$ws.Range("D56").Value = -7.984999999999999
$ws.Range("C58").Value = -13.096
$ws.Range("D61").Value = -7.634
$ws.Range("C65").Value = -12.282
$ws.Range("A68").Value = -21.526
$ws.Range("B68").Value = 5.754
$ws.Range("D69").Value = -7.311
$ws.Range("D71").Value = -7.419000000000001
$ws.Range("C73").Value = -12.646
$ws.Range("B76").Value = 5.971
$ws.Range("A77").Value = -20.196
$ws.Range("A78").Value = -19.991
$ws.Range("A79").Value = -21.054
$ws.Range("D79").Value = -7.688
$ws.Range("A80").Value = -20.534
$ws.Range("A81").Value = -21.818
$ws.Range("D81").Value = -7.650999999999999
$ws.Range("A82").Value = -22.117
$ws.Range("A84").Value = -22.013
$ws.Range("B87").Value = 4.853
$ws.Range("B88").Value = 5.16
$ws.Range("C90").Value = -13.464
$ws.Range("B92").Value = 5.401999999999999
$ws.Range("C92").Value = -10.801
$ws.Range("D92").Value = -6.731
$ws.Range("C94").Value = -10.902
$ws.Range("C95").Value = -12.001
$ws.Range("B96").Value = 6.547
$ws.Range("B98").Value = 5.89
$ws.Range("A101").Value = -21.108
$ws.Range("B101").Value = 6.234999999999999
$ws.Range("C101").Value = -12.756
$ws.Range("A102").Value = -21.358
$ws.Range("B102").Value = 6.813000000000001
